# Add a new Defect/bug entry to the "Defect" worksheet (row 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect")

$ws.Range("A9").Value = "cant set main cam after create"
$ws.Range("B9").Value = "double click thecreate could produce the bug"
$ws.Range("C9").Value = "NewCharacter"
$ws.Range("D9").Value = "Fish"
$ws.Range("E9").Value = "Fish"
$ws.Range("F9").Value = "12 Apr"

# Select E9, matching the sheet's recorded selection after the edit.
$ws.Range("E9").Select()
